# Update the "Förändrad" (Changed) date column (C) for rows 2 through 66
# from 45221 (2023-10-22) to 45224 (2023-10-25).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2:C66").Value = 45224
